$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Make Hoja2 the active sheet/tab (matches tabSelected + workbook activeTab).
$ws.Activate()

# Column widths for the new table layout (chosen so the engine's internal
# pixel-rounding of ColumnWidth lands on the stored width closest to the
# target: 24.140625 / 48.42578125 / 65.5703125 / 27.85546875 characters).
$ws.Columns.Item(1).ColumnWidth = 23.333333333333332
$ws.Columns.Item(2).ColumnWidth = 47.666666666666664
$ws.Columns.Item(3).ColumnWidth = 64.66666666666667
$ws.Columns.Item(4).ColumnWidth = 27

# Execution-path rows describing how the OCS circuit request flows through
# Util -> source (HybridSwitchSender) -> OCSSwitchSender -> HybridSwitchImpl
# down to GridSimulator. Entered in the same order the author typed them.
$ws.Range("A3").Value = "util.createOCSCircuit"
$ws.Range("B8").Value = "gridSimulator.addRequestedCircuit(ocsRoute)"
$ws.Range("B4").Value = "source.requestOCSCircuit(ocsRoute, permanent, t)"
$ws.Range("C5").Value = "requestOCSCircuit(OCSRoute ocsRoute, boolean permanent, Time time)"

# Header row (bold).
$ws.Range("A1").Value = "Util"
$ws.Range("C1").Value = "OCSSwitchSender"
$ws.Range("B1").Value = "HybridSwitchSender"
$ws.Range("D1").Value = "HybridSwitchImpl"
$ws.Range("A1:D1").Font.Bold = $true

$ws.Range("C7").Value = "GridSimulator"
$ws.Range("C7").Font.Bold = $true

$ws.Range("C9").Value = "requestedCircuits.add(route);"
$ws.Range("D6").Value = "owner.sendNow(ocsRoute.getSource(), request, time)"

# Empty, underlined marker cell.
$ws.Range("C12").Font.Underline = $true

# Print setup (paper size / orientation) for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
$null = $ws.Range("C13").Select()
